# Replace the dummy "Customer ID" value (A2) with a new dummy numeric-looking
# string, keeping it stored as text (matches the original's quote-prefixed
# text style rather than becoming a real number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'56872904"

# Replace the "Password" value (B2) with the literal text "Password" and
# select B2 as the active cell, matching the author's final saved state.
$ws.Range("B2").Value = "Password"
$ws.Range("B2").Select()
